$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# L column: numeric formula copying the population figure (F column).
# N column: string formula concatenating L (population) with M (",").
# Build L2/N2 first as standalone formulas, then fill L3:L26 / N3:N26 as a
# single fill-down operation each so the engine records them as shared
# formula groups (matching how Excel records a fill-handle drag).
$ws.Range("L2").NumberFormat = "#,##0"
$ws.Range("L2").Formula = "=F2"
$ws.Range("L3:L26").NumberFormat = "#,##0"
$ws.Range("L3:L26").Formula = "=F3"

for ($row = 2; $row -le 26; $row++) {
    $ws.Range("M$row").Value = ","
}

$ws.Range("N2").NumberFormat = "#,##0"
$ws.Range("N2").Formula = "=L2&M2"
$ws.Range("N3:N26").NumberFormat = "#,##0"
$ws.Range("N3:N26").Formula = "=L3&M3"

# Select N2:N26 (active cell N2) to match the saved view state, and scroll
# the window down so row 23 is visible.
$ws.Range("N2:N26").Select()
$excel.ActiveWindow.ScrollRow = 23

# Add print/page-setup info to the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
